$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 9.869999999999999

$ws.Range("B3").Value = 10.13
$ws.Range("D3").Value = 10.16
$ws.Range("E3").Value = 10.73

$ws.Range("C4").Value = 9.84
$ws.Range("E4").Value = 10.64
$ws.Range("F4").Value = 9.81
$ws.Range("G4").Value = 9.94

$ws.Range("C5").Value = 9.220000000000001
$ws.Range("D5").Value = 9.359999999999999
$ws.Range("F5").Value = 10.16

$ws.Range("D6").Value = 10.19
$ws.Range("E6").Value = 9.84
$ws.Range("G6").Value = 10.38

$ws.Range("D7").Value = 10.22
$ws.Range("F7").Value = 9.619999999999999
$ws.Range("I7").Value = 7.6

$ws.Range("J8").Value = 11.08

$ws.Range("G9").Value = 12.4

$ws.Range("H10").Value = 8.92
